$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "day" runtime values (rows 3-12, column B) ---
$ws.Range('B3').Value  = [double]"1.3517380000000001E-2"
$ws.Range('B4').Value  = [double]"7.6885399999999998E-3"
$ws.Range('B5').Value  = [double]"1.8429040000000001E-2"
$ws.Range('B6').Value  = [double]"3.2271800000000001E-3"
$ws.Range('B7').Value  = [double]"4.8616600000000003E-3"
$ws.Range('B8').Value  = [double]"2.1164840000000001E-2"
$ws.Range('B9').Value  = [double]"1.6294039999999999E-2"
$ws.Range('B10').Value = [double]"1.69991E-2"
$ws.Range('B11').Value = [double]"1.8328480000000001E-2"
$ws.Range('B12').Value = [double]"0.12447448"

# --- Append new rows for day 11, 12 and 13 ---
$ws.Range('A13').Value = 11
$ws.Range('B13').Value = [double]"0.10642092"

$ws.Range('A14').Value = 12
$ws.Range('B14').Value = [double]"9.5150639999999995E-2"

$ws.Range('A15').Value = 13
$ws.Range('B15').Value = [double]"2.3314400000000002E-3"

# --- Extend the bar chart's series references to cover the new data ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(RuntimesChart!`$B`$2,RuntimesChart!`$A`$3:`$A`$20,RuntimesChart!`$B`$3:`$B`$20,1)"

# --- Match the selection left active in the sheet after the edits ---
$ws.Range("D20").Select()
